$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$targetRows = @(2, 4, 5, 7, 8, 9, 10, 12, 13, 15, 16)
foreach ($r in $targetRows) {
    $cell = $t.Cell($r, 5)
    $para = $cell.Range.Paragraphs.Item(1)
    $para.Range.Style = "Normal"
}

Write-Output "done"
